$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "48.373.98"
$ws.Range("E2").Value = "  +2.07%  "
$ws.Range("D3").Value = "2.524.39"
$ws.Range("E3").Value = "  +0.74%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "324.18"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -0.03%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "109.50"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -0.44%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.528"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  +0.49%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.999"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -0.02%  "
$ws.Range("E9").Value = "  +4.46%  "
$ws.Range("E10").Value = "  +1.99%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "19.88"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +7.09%  "
$ws.Range("E12").Value = "  +0.46%  "
$ws.Range("E13").Value = "  +0.92%  "
$ws.Range("E14").Value = "  +0.28%  "
$ws.Range("D15").Value = "2.920.47"
$ws.Range("E15").Value = "  +0.71%  "
$ws.Range("D16").Value = "2.524.08"
$ws.Range("E16").Value = "  +0.63%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.857"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  -0.69%  "
$ws.Range("D18").Value = "48.235.69"
$ws.Range("E18").Value = "  +1.92%  "
$ws.Range("E19").Value = "  +3.59%  "
$ws.Range("E20").Value = "  -0.40%  "
$ws.Range("E21").Value = "  +0.13%  "
$ws.Range("E22").Value = "  +2.66%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "72.81"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +3.12%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "273.48"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +9.71%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.59"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -0.65%  "
$ws.Range("B26").Value = "Dai"
$ws.Range("C26").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.00"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +0.14%  "
$ws.Range("B27").Value = "EthereumClassic"
$ws.Range("C27").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "26.14"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +0.01%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.19"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +1.45%  "
$ws.Range("B29").Value = "Toncoin"
$ws.Range("C29").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.27"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -1.05%  "
$ws.Range("B30").Value = "Kaspa"
$ws.Range("C30").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.146"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +5.76%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "35.36"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -0.77%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "49.79"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -0.99%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "20.01"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -0.54%  "
$ws.Range("E34").Value = "  -0.11%  "
$ws.Range("E35").Value = "  -0.02%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0791"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -1.23%  "
$ws.Range("E37").Value = "  -0.77%  "
$ws.Range("E38").Value = "  -0.23%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.99"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -0.57%  "
$ws.Range("E40").Value = "  +0.04%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "22.23"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +3.31%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "119.24"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -2.20%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.18"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -1.77%  "
$ws.Range("E44").Value = "  +0.36%  "
$ws.Range("D45").Value = "2.003.38"
$ws.Range("E45").Value = "  -0.14%  "
$ws.Range("E46").Value = "  +1.12%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.88"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +5.49%  "
$ws.Range("E48").Value = "  -2.46%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "9.10"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +0.47%  "
$ws.Range("E50").Value = "  +0.64%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "81.11"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +3.40%  "
